$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers (1-based) that receive updated values: B,C,D,E,F,H,J,K,L,M
$cols = @(2,3,4,5,6,8,10,11,12,13)

# Each inner array holds the new values for columns B,C,D,E,F,H,J,K,L,M, row by row (rows 2..25)
$data = @(
    @(21.00161372784111,2.971355352184507,4.486033468703581,10.53686874755586,71.63293216806579,7.344005520526261,11.31345065817433,17.29396089744493,10.81972464555595,18.70543065279006),
    @(21.03016882288739,2.819988934432347,4.35641589113319,10.55991303917255,70.82504543664341,7.344005520526261,11.30218204383506,17.33297963058585,10.85970154143519,18.76925023135851),
    @(21.05400610868024,2.776838938734952,4.273843329462611,10.5749019525427,70.32644243018953,7.344005520526261,11.2953799438293,17.36201502296641,10.88579955826682,18.8120863576133),
    @(21.06530247405447,2.765640278106269,4.239464926892398,10.58122186789469,70.12274315558327,7.344005520526261,11.29263745090001,17.37512086623438,10.8968258944963,18.83046018043168),
    @(21.06727370778163,2.763815866735772,4.233713036038255,10.58228409795163,70.08889171700415,7.344005520526261,11.29218383947581,17.3773739084959,10.89868046253703,18.83356656266563),
    @(21.05415205235946,2.77668557071437,4.2733826135656,10.57498632660982,70.3236971845467,7.344005520526261,11.29534283845366,17.36218662039015,10.8859466782143,18.81233043782427),
    @(21.01014960125462,2.920193009801759,4.44197362530051,10.54464056480838,71.35498077092916,7.344005520526261,11.30954099955477,17.30635890308243,10.83318709888837,18.72667781879587),
    @(20.97397554205977,3.270192183098579,4.748030900005877,10.49176393706369,73.35140493925323,7.344005520526261,11.33831069465593,17.2372936729259,10.74199967621597,18.58768980651886),
    @(20.97802747226087,3.502869419750877,4.956923533248776,10.45691557526495,74.79490540097447,7.344005520526261,11.36001433475518,17.21132990371358,10.68243069661479,18.50325401577473),
    @(20.98652062791072,3.603360837487464,5.048310776973752,10.44192171643683,75.44497729968298,7.344005520526261,11.37001109650154,17.20491903955061,10.65693205546888,18.46868341311775),
    @(20.99069138652688,3.640642027955794,5.082380787547222,10.43636673720735,75.69006696623659,7.344005520526261,11.37381433758496,17.203268723157,10.64750551948021,18.45614477408999),
    @(20.98975071984358,3.632647255573057,5.075067247727012,10.4375576466923,75.63733255273175,7.344005520526261,11.37299445674565,17.20358956905601,10.64952550904313,18.45882061645368),
    @(20.98684464021899,3.60644348178372,5.051124556063645,10.44146224605133,75.46516333225583,7.344005520526261,11.37032363761716,17.20476768537084,10.65615193862611,18.46764077691089),
    @(20.98518883264353,3.590292219692476,5.036388754943255,10.44386990855157,75.3595604058101,7.344005520526261,11.36868998369266,17.20559056246327,10.66024064997909,18.47311534153976),
    @(20.97760611640862,3.496194059895134,4.950876783272356,10.45791268722462,74.75227786946859,7.344005520526261,11.35936350316795,17.21185761499067,10.6841292141987,18.50559058117024),
    @(20.97465667351926,3.437094115653057,4.897476125670279,10.46674699587571,74.37794887930464,7.344005520526261,11.35367361675758,17.21708610615866,10.69919321303736,18.52649678337164),
    @(20.9735863108687,3.402597471053258,4.866419634195732,10.47190913161481,74.16203277682629,7.344005520526261,11.35041259132657,17.22060167348851,10.70800823277239,18.53888286557612),
    @(20.97333146977966,3.390831053862734,4.855846150616213,10.47367085170861,74.08882624669016,7.344005520526261,11.34931046732141,17.22187924802586,10.71101873746321,18.54313863737916),
    @(20.97490585928113,3.443437582778181,4.903196195891313,10.46579820296721,74.41786113888196,7.344005520526261,11.35427811014497,17.21647691218839,10.69757404318,18.52423387762432),
    @(20.98767233688919,3.614161153217409,5.058171758139438,10.440312040919,75.51576391031728,7.344005520526261,11.37110764283086,17.20440054364109,10.65419938004859,18.46503508463333),
    @(21.00157860032989,3.721234309968454,5.15632594146116,10.42437125041275,76.22696178893776,7.344005520526261,11.38220973984454,17.20103884649521,10.62718740266663,18.42956577660722),
    @(20.99364837102727,3.664500141267599,5.10422965168154,10.4328138532237,75.84800543010296,7.344005520526261,11.37627495043761,17.20241834680792,10.6414822226588,18.44820162479989),
    @(20.97479125449514,3.440571318623078,4.900611259250299,10.46622689311935,74.3998190057646,7.344005520526261,11.35400478702377,17.21675074136784,10.69830558847596,18.52525579469699),
    @(20.97838506227195,3.17978378673307,4.667972888190985,10.50536295188873,72.81486936727323,7.344005520526261,11.33042828990341,17.25163462422744,10.7653601280401,18.62218609169955)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowValues[$j]
    }
}

Write-Output "Updated $($data.Count) rows"